# Applies the "find bks, comparison table and statistical analysis working"
# update: Resumen!B2/C2 (best-known zone + value), Solucion!B column
# (Salida codes reshuffled per Pedido row), and Metricas!B2:B5 (Tiempo
# values recomputed per zone).

$wb = $excel.ActiveWorkbook

$wsResumen  = $wb.Worksheets.Item("Resumen")
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

# --- Resumen: winning zone moved from Z1 to Z3, with its matching time ---
$wsResumen.Range("B2").Value = "Z3"
$wsResumen.Range("C2").Value = 601.1304735647099

# --- Solucion: Pedido -> Salida assignments reshuffled ---
$wsSolucion.Range("B2").Value = "S031"
$wsSolucion.Range("B4").Value = "S067"
$wsSolucion.Range("B5").Value = "S065"
$wsSolucion.Range("B7").Value = "S073"
$wsSolucion.Range("B8").Value = "S077"
$wsSolucion.Range("B11").Value = "S041"
$wsSolucion.Range("B12").Value = "S074"
$wsSolucion.Range("B14").Value = "S046"
$wsSolucion.Range("B15").Value = "S078"
$wsSolucion.Range("B16").Value = "S039"
$wsSolucion.Range("B17").Value = "S070"
$wsSolucion.Range("B18").Value = "S033"
$wsSolucion.Range("B19").Value = "S075"
$wsSolucion.Range("B20").Value = "S002"
$wsSolucion.Range("B23").Value = "S005"
$wsSolucion.Range("B24").Value = "S049"
$wsSolucion.Range("B25").Value = "S079"
$wsSolucion.Range("B26").Value = "S038"
$wsSolucion.Range("B27").Value = "S068"
$wsSolucion.Range("B28").Value = "S003"
$wsSolucion.Range("B29").Value = "S040"
$wsSolucion.Range("B30").Value = "S009"
$wsSolucion.Range("B31").Value = "S008"
$wsSolucion.Range("B32").Value = "S035"
$wsSolucion.Range("B33").Value = "S028"
$wsSolucion.Range("B34").Value = "S027"
$wsSolucion.Range("B35").Value = "S042"
$wsSolucion.Range("B36").Value = "S044"
$wsSolucion.Range("B37").Value = "S032"
$wsSolucion.Range("B38").Value = "S037"
$wsSolucion.Range("B39").Value = "S030"
$wsSolucion.Range("B40").Value = "S072"
$wsSolucion.Range("B42").Value = "S053"
$wsSolucion.Range("B43").Value = "S012"
$wsSolucion.Range("B44").Value = "S001"
$wsSolucion.Range("B45").Value = "S076"
$wsSolucion.Range("B46").Value = "S043"
$wsSolucion.Range("B47").Value = "S045"
$wsSolucion.Range("B48").Value = "S047"
$wsSolucion.Range("B49").Value = "S010"
$wsSolucion.Range("B50").Value = "S056"
$wsSolucion.Range("B51").Value = "S036"
$wsSolucion.Range("B52").Value = "S054"
$wsSolucion.Range("B53").Value = "S004"
$wsSolucion.Range("B54").Value = "S048"
$wsSolucion.Range("B55").Value = "S007"
$wsSolucion.Range("B56").Value = "S014"
$wsSolucion.Range("B57").Value = "S013"
$wsSolucion.Range("B58").Value = "S006"
$wsSolucion.Range("B59").Value = "S055"
$wsSolucion.Range("B60").Value = "S050"
$wsSolucion.Range("B61").Value = "S052"
$wsSolucion.Range("B62").Value = "S051"
$wsSolucion.Range("B63").Value = "S015"
$wsSolucion.Range("B64").Value = "S016"
$wsSolucion.Range("B65").Value = "S011"
$wsSolucion.Range("B66").Value = "S021"
$wsSolucion.Range("B68").Value = "S058"
$wsSolucion.Range("B69").Value = "S061"
$wsSolucion.Range("B70").Value = "S018"
$wsSolucion.Range("B71").Value = "S062"
$wsSolucion.Range("B72").Value = "S022"
$wsSolucion.Range("B73").Value = "S017"
$wsSolucion.Range("B74").Value = "S019"
$wsSolucion.Range("B75").Value = "S059"
$wsSolucion.Range("B77").Value = "S063"
$wsSolucion.Range("B79").Value = "S060"
$wsSolucion.Range("B80").Value = "S024"

# --- Metricas: recomputed Tiempo per zone ---
$wsMetricas.Range("B2").Value = 601.0939290734133
$wsMetricas.Range("B3").Value = 518.4691317980323
$wsMetricas.Range("B4").Value = 601.1304735647099
$wsMetricas.Range("B5").Value = 554.3820953616606
